$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row64 - Forged from the Void / Void Glue
$ws.Range("H64").Value = 4428.5713
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4428.5713
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4428.5713
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4924.5713

# ALC!row67 - Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 4428.5713
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4428.5713
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4428.5713
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6144.5713

# ALC!row70 - Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 2574.3125
$ws.Range("I70").Value = 1928.6923
$ws.Range("J70").Value = 5372
$ws.Range("K70").Value = 5786.0769
$ws.Range("L70").Value = 16116
$ws.Range("M70").Value = -5516.0769
$ws.Range("N70").Value = -16656

# ALC!row73 - Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 2574.3125
$ws.Range("I73").Value = 1928.6923
$ws.Range("J73").Value = 5372
$ws.Range("K73").Value = 5786.0769
$ws.Range("L73").Value = 16116
$ws.Range("M73").Value = -4850.0769
$ws.Range("N73").Value = -17988

# ALC!row76 - Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 7624
$ws.Range("I76").Value = 7624
$ws.Range("K76").Value = 7624
$ws.Range("M76").Value = -7309

# ALC!row79 - The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 7624
$ws.Range("I79").Value = 7624
$ws.Range("K79").Value = 7624
$ws.Range("M79").Value = -6532

# ALC!row92 - Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 1090.3
$ws.Range("I92").Value = 1573.4615
$ws.Range("J92").Value = 193
$ws.Range("K92").Value = 1573.4615
$ws.Range("L92").Value = 193
$ws.Range("M92").Value = -325.4614999999999
$ws.Range("N92").Value = -2689

# ALC!row113 - Amaro Kart / Starch Glue
$ws.Range("H113").Value = 11010.471
$ws.Range("I113").Value = 11620.9
$ws.Range("J113").Value = 10138.429
$ws.Range("K113").Value = 11620.9
$ws.Range("L113").Value = 10138.429
$ws.Range("M113").Value = -8366.9
$ws.Range("N113").Value = -16646.429

# ALC!row132 - Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 11069.594
$ws.Range("I132").Value = 9360.76
$ws.Range("K132").Value = 28082.28
$ws.Range("M132").Value = -25552.28

# ALC!row137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 11883.4375
$ws.Range("I137").Value = 1858
$ws.Range("J137").Value = 18742.947
$ws.Range("K137").Value = 5574
$ws.Range("L137").Value = 56228.841
$ws.Range("M137").Value = -3024
$ws.Range("N137").Value = -61328.841

$ws = $wb.Worksheets.Item("ARM")
# ARM!row41 - Skillet Scandal / White Skillet
$ws.Range("H41").Value = 5538
$ws.Range("I41").Value = 363.33334
$ws.Range("J41").Value = 21062
$ws.Range("K41").Value = 363.33334
$ws.Range("L41").Value = 21062
$ws.Range("M41").Value = 50.66665999999998
$ws.Range("N41").Value = -21890

# ARM!row61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 15539.482
$ws.Range("I61").Value = 4153.875
$ws.Range("J61").Value = 19876.857
$ws.Range("K61").Value = 4153.875
$ws.Range("L61").Value = 19876.857
$ws.Range("M61").Value = -3941.875
$ws.Range("N61").Value = -20300.857

# ARM!row63 - Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 2441.1428
$ws.Range("I63").Value = 2264.6667
$ws.Range("K63").Value = 2264.6667
$ws.Range("M63").Value = -1578.6667

# ARM!row66 - A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 2441.1428
$ws.Range("I66").Value = 2264.6667
$ws.Range("K66").Value = 11323.3335
$ws.Range("M66").Value = -7891.333500000001

# ARM!row122 - Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 3255.7856
$ws.Range("I122").Value = 2028.1
$ws.Range("J122").Value = 6325
$ws.Range("K122").Value = 6084.299999999999
$ws.Range("L122").Value = 18975
$ws.Range("M122").Value = -3634.299999999999
$ws.Range("N122").Value = -23875

# ARM!row136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 15539.482
$ws.Range("I136").Value = 4153.875
$ws.Range("J136").Value = 19876.857
$ws.Range("K136").Value = 12461.625
$ws.Range("L136").Value = 59630.571
$ws.Range("M136").Value = -9911.625
$ws.Range("N136").Value = -64730.571

$ws = $wb.Worksheets.Item("BSM")
# BSM!row106 - Fire for Hire / Molybdenum Rimfire
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524

$ws = $wb.Worksheets.Item("CRP")
# CRP!row62 - Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 6000
$ws.Range("J62").Value = 6000
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7248

# CRP!row65 - The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 6000
$ws.Range("J65").Value = 6000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240

# CRP!row122 - Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 2397.75
$ws.Range("J122").Value = 2400
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100

# CRP!row134 - Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 40008420
$ws.Range("I134").Value = 1890.5
$ws.Range("K134").Value = 5671.5
$ws.Range("M134").Value = -3136.5

$ws = $wb.Worksheets.Item("CUL")
# CUL!row5 - What a Sap / Maple Syrup
$ws.Range("H5").Value = 4066374.2
$ws.Range("I5").Value = 1600.4
$ws.Range("J5").Value = 24390244
$ws.Range("K5").Value = 4801.200000000001
$ws.Range("L5").Value = 73170732
$ws.Range("M5").Value = -4689.200000000001
$ws.Range("N5").Value = -73170956

# CUL!row17 - Chew the Fat / Grilled Dodo
$ws.Range("H17").Value = 227
$ws.Range("I17").Value = 345
$ws.Range("J17").Value = 138.5
$ws.Range("K17").Value = 1035
$ws.Range("L17").Value = 415.5
$ws.Range("M17").Value = -866
$ws.Range("N17").Value = -753.5

# CUL!row122 - Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 12196.4
$ws.Range("J122").Value = 25742.5
$ws.Range("L122").Value = 231682.5
$ws.Range("N122").Value = -236582.5

# CUL!row131 - The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1495.7878
$ws.Range("J131").Value = 1498.7938
$ws.Range("L131").Value = 4496.3814
$ws.Range("N131").Value = -14576.3814

# CUL!row132 - More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 2015.4445
$ws.Range("I132").Value = 1892.375
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 17031.375
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -14501.375
$ws.Range("N132").Value = -32060

# CUL!row135 - Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 4066374.2
$ws.Range("I135").Value = 1600.4
$ws.Range("J135").Value = 24390244
$ws.Range("K135").Value = 14403.6
$ws.Range("L135").Value = 219512196
$ws.Range("M135").Value = -11868.6
$ws.Range("N135").Value = -219517266

$ws = $wb.Worksheets.Item("GSM")
# GSM!row122 - Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3179.8572
$ws.Range("I122").Value = 3179.8572
$ws.Range("K122").Value = 9539.571599999999
$ws.Range("M122").Value = -7089.571599999999

$ws = $wb.Worksheets.Item("LTW")
# LTW!row122 - Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 5772.409
$ws.Range("I122").Value = 5166.2
$ws.Range("K122").Value = 15498.6
$ws.Range("M122").Value = -13048.6

# LTW!row136 - Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 14178.579
$ws.Range("I136").Value = 19499.5
$ws.Range("J136").Value = 11722.77
$ws.Range("K136").Value = 58498.5
$ws.Range("L136").Value = 35168.31
$ws.Range("M136").Value = -55948.5
$ws.Range("N136").Value = -40268.31

$ws = $wb.Worksheets.Item("WVR")
# WVR!row122 - Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 2386.3416
$ws.Range("I122").Value = 1330.6451
$ws.Range("J122").Value = 5659
$ws.Range("K122").Value = 3991.9353
$ws.Range("L122").Value = 16977
$ws.Range("M122").Value = -1541.9353
$ws.Range("N122").Value = -21877
